$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3: 85.59999999999999 -> 84.59999999999999
$ws.Range("D3").Value = 84.59999999999999

# Update E12: blank numeric cell -> text date string "2022-08-29"
# Enter it as a text-producing formula first (Excel won't re-parse a formula
# result as a date the way it would a typed literal), then convert it to a
# plain value in place. This keeps the cell's existing style/format (s="4",
# matching the other Driver Vintage cells like E13/E14) instead of Excel
# silently switching the cell to a date or "quoted text" style.
$ws.Range("E12").Formula = "=""2022-08-29"""
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
